$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 68 and 69 swap their entire record content (everything except the
# shared location/date/observer columns, which are identical for both rows
# already).
# ---------------------------------------------------------------------------

# New content for row 68 (previously held by row 69)
$ws.Range("A68").Value = 130837316
$ws.Range("B68").Value = 83223
$ws.Range("E68").Value = 6440
$ws.Range("F68").Value = "Vitgrynig nållav"
$ws.Range("G68").Value = "Chaenotheca subroscida"
$ws.Range("H68").Value = "(Eitner) Zahlbr."
$ws.Range("I68").ClearContents()
$ws.Range("J68").ClearContents()
$ws.Range("Q68").Value = 445777
$ws.Range("R68").Value = 7026331
$ws.Range("S68").Value = 6
$ws.Range("Z68").Value = "10:36"
$ws.Range("AB68").Value = "10:36"
$ws.Range("AC68").Value = "Rikligt på bark på stam av levande gammal gran i gammal granskog"

# New content for row 69 (previously held by row 68)
$ws.Range("A69").Value = 130839110
$ws.Range("B69").Value = 79243
$ws.Range("E69").Value = 6425
$ws.Range("F69").Value = "Garnlav"
$ws.Range("G69").Value = "Alectoria sarmentosa"
$ws.Range("H69").Value = "(Ach.) Ach."
$ws.Range("I69").NumberFormat = "@"
$ws.Range("I69").Value = "75"
$ws.Range("J69").Value = "bålar"
$ws.Range("Q69").Value = 445730
$ws.Range("R69").Value = 7026205
$ws.Range("S69").Value = 5
$ws.Range("Z69").Value = "12:29"
$ws.Range("AB69").Value = "12:29"
$ws.Range("AC69").Value = "På gammal klen död gran i gammal granskog"

# ---------------------------------------------------------------------------
# Rows 77, 78, 79 rotate their content: new77 = old78, new78 = old79,
# new79 = old77.
# ---------------------------------------------------------------------------

# New content for row 77 (previously held by row 78)
$ws.Range("A77").Value = 130839413
$ws.Range("B77").Value = 78255
$ws.Range("D77").Value = "NT"
$ws.Range("E77").Value = 228579
$ws.Range("F77").Value = "Liten svartspik"
$ws.Range("G77").Value = "Chaenothecopsis nana"
$ws.Range("H77").Value = "Tibell"
$ws.Range("Q77").Value = 445781
$ws.Range("R77").Value = 7026373
$ws.Range("S77").Value = 7
$ws.Range("Z77").Value = "12:54"
$ws.Range("AB77").Value = "12:54"
$ws.Range("AC77").Value = "På bark på stam av levande gammal gran i gles gammal granskog"

# New content for row 78 (previously held by row 79)
$ws.Range("A78").Value = 130837541
$ws.Range("B78").Value = 75221
$ws.Range("D78").Value = "LC"
$ws.Range("E78").Value = 6428
$ws.Range("F78").Value = "Rostfläck"
$ws.Range("G78").Value = "Arthonia vinosa"
$ws.Range("H78").Value = "Leight."
$ws.Range("Q78").Value = 445740
$ws.Range("R78").Value = 7026322
$ws.Range("S78").Value = 8
$ws.Range("Z78").Value = "10:52"
$ws.Range("AB78").Value = "10:52"
$ws.Range("AC78").Value = "På tunna kvistar vid basen på gammal levande gran"

# New content for row 79 (previously held by row 77)
$ws.Range("A79").Value = 130837733
$ws.Range("B79").Value = 79243
$ws.Range("D79").Value = "NT"
$ws.Range("E79").Value = 6425
$ws.Range("F79").Value = "Garnlav"
$ws.Range("G79").Value = "Alectoria sarmentosa"
$ws.Range("H79").Value = "(Ach.) Ach."
$ws.Range("Q79").Value = 445720
$ws.Range("R79").Value = 7026343
$ws.Range("S79").Value = 10
$ws.Range("Z79").Value = "10:59"
$ws.Range("AB79").Value = "10:59"
$ws.Range("AC79").Value = "På gammal död gran i gammal granskog"

# ---------------------------------------------------------------------------
# New row 81: a new observation record (Lavskrika / Perisoreus infaustus).
# ---------------------------------------------------------------------------

$ws.Range("A81").Value = 131074382
$ws.Range("B81").Value = 57988
$ws.Range("D81").Value = "LC"
$ws.Range("E81").Value = 103031
$ws.Range("F81").Value = "Lavskrika"
$ws.Range("G81").Value = "Perisoreus infaustus"
$ws.Range("H81").Value = "(Linnaeus, 1758)"
$ws.Range("I81").NumberFormat = "@"
$ws.Range("I81").Value = "1"
$ws.Range("M81").Value = "lockläte, övriga läten"
$ws.Range("P81").Value = "Svedaun, Svedaun, Jmt"
$ws.Range("Q81").Value = 445791
$ws.Range("R81").Value = 7026340
$ws.Range("S81").Value = 15
$ws.Range("T81").Value = "Jämtland"
$ws.Range("U81").Value = "Krokom"
$ws.Range("V81").Value = "Jämtland"
$ws.Range("W81").Value = "Alsen"
$ws.Range("Y81").NumberFormat = "@"
$ws.Range("Y81").Value = "2026-02-08"
$ws.Range("Z81").Value = "12:12"
$ws.Range("AA81").NumberFormat = "@"
$ws.Range("AA81").Value = "2026-02-08"
$ws.Range("AB81").Value = "12:12"
$ws.Range("AC81").Value = "Lockade i ca 3 minuter på en grantopp"
$ws.Range("AD81").Value = $false
$ws.Range("AE81").Value = $false
$ws.Range("AG81").Value = $false
$ws.Range("AW81").Value = "Ludvig Nordin"
$ws.Range("AX81").Value = "Ludvig Nordin"
